# Fill in the previously-blank placeholder rows (608-614) with the
# electricity/price data for Jul 2023 - Jan 2024, and move the current
# selection down to the still-empty rows below (615:AC626), matching the
# author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A=Month, B..M = electricity figures, N..P = price figures
$rowData = @{
    608 = @(45108, 412.235, 1.566, 12.102, 425.902, 3.217, 1.963, 1.254, 28.543, 386.432, 12.181, 398.612, 24.07, 27.43, 26.6,  25.78)
    609 = @(45139, 410.087, 1.542, 12.413, 424.042, 3.001, 1.606, 1.395, 21.1,   391.9,   12.436, 404.336, 23.39, 24.92, 22.46, 23.97)
    610 = @(45170, 345.956, 1.427, 11.664, 359.047, 2.356, 2.25,  0.106, 1.357,  346.129, 11.667, 357.795, 19.68, 21.16, 19.44, 20.34)
    611 = @(45200, 316.802, 1.364, 11.33,  329.497, 2.211, 2.006, 0.205, 10.514, 307.874, 11.314, 319.188, 24.02, 25.21, 25,    24.66)
    612 = @(45231, 308.934, 1.393, 11.776, 322.103, 2.469, 1.829, 0.64,  17.52,  293.487, 11.737, 305.224, 21.44, 21.92, 20.19, 21.55)
    613 = @(45261, 332.392, 1.462, 12.534, 346.387, 3.202, 1.863, 1.339, 24.294, 310.959, 12.473, 323.432, 19.44, 20.16, 19.14, 19.76)
    614 = @(45292, 365.625, 1.481, 12.693, 379.799, 3.613, 1.876, 1.737, 27.894, 341.01,  12.632, 353.643, 21.96, 22.66, 21.8,  22.3)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Move the selection to the still-empty trailing block, as saved by the author.
$ws.Range("A615:AC626").Select()
